$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Hours spent added for existing "TankHistory page" row
$ws.Range("H3").Value = 10

# New row 4: Yard page
$ws.Range("A4").Value = "Yard page"
$ws.Range("B4").Value = 8
$ws.Range("G4").Value = 7
$ws.Rows.Item(4).RowHeight = 40.5

# New row 5: Invoices page
$ws.Range("A5").Value = "Invoices page"
$ws.Range("B5").Value = 8
$ws.Range("G5").Value = 6
$ws.Rows.Item(5).RowHeight = 40.5

# Update the view/selection to match the saved state
$ws.Range("G5").Select()
$excel.ActiveWindow.ScrollColumn = 7
